$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

Set-TextValue "D2" "48.124.17"
Set-TextValue "E2" "  -0.12%  "
Set-TextValue "D3" "2.497.08"
Set-TextValue "E3" "  -1.39%  "
Set-TextValue "E4" "  -0.14%  "
Set-TextValue "D5" "318.08"
Set-TextValue "E5" "  -1.90%  "
Set-TextValue "D6" "105.89"
Set-TextValue "E6" "  -2.80%  "
Set-TextValue "D7" "0.520"
Set-TextValue "E7" "  -1.55%  "
Set-TextValue "E8" "  -0.09%  "
Set-TextValue "D9" "0.539"
Set-TextValue "E9" "  -3.60%  "
Set-TextValue "E10" "  -4.55%  "
Set-TextValue "D11" "20.28"
Set-TextValue "E11" "  -0.95%  "
Set-TextValue "E12" "  -3.20%  "
Set-TextValue "E13" "  -0.15%  "
Set-TextValue "E14" "  -2.80%  "
Set-TextValue "D15" "2.888.75"
Set-TextValue "E15" "  -1.38%  "
Set-TextValue "D16" "2.497.26"
Set-TextValue "E16" "  -1.40%  "
Set-TextValue "D17" "0.829"
Set-TextValue "E17" "  -3.66%  "
Set-TextValue "D18" "48.013.73"
Set-TextValue "E18" "  -0.09%  "
Set-TextValue "D19" "2.99"
Set-TextValue "E19" "  +11.41%  "
Set-TextValue "D20" "12.81"
Set-TextValue "E20" "  -3.51%  "
Set-TextValue "D21" "6.58"
Set-TextValue "E21" "  -1.05%  "
Set-TextValue "D22" "0.0₃0931"
Set-TextValue "E22" "  -2.22%  "
Set-TextValue "D23" "71.15"
Set-TextValue "E23" "  -1.66%  "
Set-TextValue "D24" "267.83"
Set-TextValue "E24" "  -0.85%  "
Set-TextValue "E25" "  -2.80%  "
Set-TextValue "E26" "  +0.17%  "
Set-TextValue "D27" "25.77"
Set-TextValue "E27" "  -1.88%  "
Set-TextValue "D28" "2.28"
Set-TextValue "E28" "  +3.26%  "
Set-TextValue "D29" "9.76"
Set-TextValue "E29" "  -4.04%  "
Set-TextValue "E30" "  -4.01%  "
Set-TextValue "D31" "34.61"
Set-TextValue "D32" "49.37"
Set-TextValue "E32" "  -0.85%  "
Set-TextValue "D34" "19.12"
Set-TextValue "E34" "  -4.06%  "
Set-TextValue "E35" "  -2.48%  "
Set-TextValue "E36" "  -2.80%  "
Set-TextValue "E37" "  -2.68%  "
Set-TextValue "E38" "  -3.37%  "
Set-TextValue "E39" "  -4.53%  "
Set-TextValue "D40" "122.91"
Set-TextValue "E40" "  +3.05%  "
Set-TextValue "D41" "22.35"
Set-TextValue "E41" "  -0.25%  "
Set-TextValue "E42" "  -1.90%  "
Set-TextValue "E43" "  +1.35%  "
Set-TextValue "E44" "  +0.29%  "
Set-TextValue "D45" "2.000.91"
Set-TextValue "E45" "  -0.71%  "
Set-TextValue "D46" "3.15"
Set-TextValue "E46" "  -0.23%  "
Set-TextValue "D47" "1.91"
Set-TextValue "E47" "  +1.64%  "
Set-TextValue "E48" "  -2.73%  "
Set-TextValue "D49" "8.94"
Set-TextValue "E49" "  -2.62%  "
Set-TextValue "D50" "5.19"
Set-TextValue "E50" "  -1.53%  "
Set-TextValue "D51" "78.80"
Set-TextValue "E51" "  -1.34%  "
